$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = '''62.152.88'
$ws.Range("E2").Formula = '''  -2.98%  '
$ws.Range("D3").Formula = '''3.183.81'
$ws.Range("E3").Formula = '''  -4.40%  '
$ws.Range("E4").Formula = '''  +0.01%  '
$ws.Range("D5").Formula = '''585.59'
$ws.Range("E5").Formula = '''  -2.85%  '
$ws.Range("D6").Formula = '''135.07'
$ws.Range("E6").Formula = '''  -6.34%  '
$ws.Range("E7").Formula = '''  -0.09%  '
$ws.Range("D8").Formula = '''3.182.79'
$ws.Range("E8").Formula = '''  -4.41%  '
$ws.Range("D9").Formula = '''0.503'
$ws.Range("E9").Formula = '''  -4.22%  '
$ws.Range("E10").Formula = '''  -6.15%  '
$ws.Range("E11").Formula = '''  -6.01%  '
$ws.Range("D12").Formula = '''0.451'
$ws.Range("E12").Formula = '''  -5.47%  '
$ws.Range("E13").Formula = '''  -6.89%  '
$ws.Range("D14").Formula = '''33.16'
$ws.Range("E14").Formula = '''  -5.47%  '
$ws.Range("D15").Formula = '''3.705.63'
$ws.Range("E15").Formula = '''  -4.46%  '
$ws.Range("D16").Formula = '''0.119'
$ws.Range("E16").Formula = '''  -1.52%  '
$ws.Range("D17").Formula = '''3.178.99'
$ws.Range("E17").Formula = '''  -4.35%  '
$ws.Range("D18").Formula = '''62.243.73'
$ws.Range("E18").Formula = '''  -2.98%  '
$ws.Range("D19").Formula = '''6.58'
$ws.Range("E19").Formula = '''  -5.19%  '
$ws.Range("D20").Formula = '''455.52'
$ws.Range("E20").Formula = '''  -6.08%  '
$ws.Range("D21").Formula = '''13.99'
$ws.Range("E21").Formula = '''  -2.74%  '
$ws.Range("D22").Formula = '''0.704'
$ws.Range("E22").Formula = '''  -5.16%  '
$ws.Range("D23").Formula = '''7.61'
$ws.Range("E23").Formula = '''  -5.54%  '
$ws.Range("D24").Formula = '''13.38'
$ws.Range("E24").Formula = '''  -3.45%  '
$ws.Range("D25").Formula = '''82.44'
$ws.Range("E25").Formula = '''  -3.16%  '
$ws.Range("E26").Formula = '''  -0.21%  '
$ws.Range("D27").Formula = '''0.999'
$ws.Range("E27").Formula = '''  +0.02%  '
$ws.Range("E28").Formula = '''  -4.28%  '
$ws.Range("D29").Formula = '''6.89'
$ws.Range("E29").Formula = '''  -5.46%  '
$ws.Range("D30").Formula = '''7.81'
$ws.Range("E30").Formula = '''  -6.77%  '
$ws.Range("D31").Formula = '''2.02'
$ws.Range("E31").Formula = '''  -7.28%  '
$ws.Range("D32").Formula = '''27.27'
$ws.Range("E32").Formula = '''  -8.66%  '
$ws.Range("E33").Formula = '''  -4.50%  '
$ws.Range("E34").Formula = '''  -8.12%  '
$ws.Range("E35").Formula = '''  -6.32%  '
$ws.Range("D36").Formula = '''5.79'
$ws.Range("E36").Formula = '''  -4.16%  '
$ws.Range("D37").Formula = '''51.03'
$ws.Range("E37").Formula = '''  -4.49%  '
$ws.Range("D38").Formula = '''0.0₃0688'
$ws.Range("E38").Formula = '''  -10.06%  '
$ws.Range("D39").Formula = '''0.0383'
$ws.Range("E39").Formula = '''  -5.27%  '
$ws.Range("B40").Value = 'Bittensor'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D40").Formula = '''408.90'
$ws.Range("E40").Formula = '''  -6.53%  '
$ws.Range("B41").Value = 'Maker'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D41").Formula = '''2.947.03'
$ws.Range("E41").Formula = '''  -3.69%  '
$ws.Range("B42").Value = 'Kaspa'
$ws.Range("C42").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D42").Formula = '''0.113'
$ws.Range("E42").Formula = '''  +1.03%  '
$ws.Range("B43").Value = 'Cosmos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D43").Formula = '''8.02'
$ws.Range("E43").Formula = '''  -5.45%  '
$ws.Range("D44").Formula = '''2.61'
$ws.Range("E44").Formula = '''  -7.18%  '
$ws.Range("E45").Formula = '''  -7.70%  '
$ws.Range("D46").Formula = '''2.14'
$ws.Range("E46").Formula = '''  -4.90%  '
$ws.Range("E47").Formula = '''  -0.07%  '
$ws.Range("D48").Formula = '''35.81'
$ws.Range("E48").Formula = '''  -1.29%  '
$ws.Range("D49").Formula = '''25.44'
$ws.Range("E49").Formula = '''  -4.87%  '
$ws.Range("D50").Formula = '''122.81'
$ws.Range("E50").Formula = '''  -0.77%  '
$ws.Range("E51").Formula = '''  -4.46%  '
